$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "update ds rule import data": add a new "percentage" column (E) to the
# customer_code / customer_name / start_date / end_date import template,
# with a sample value of 30 on the existing data row.
$ws.Range("E1").Value = "percentage"
$ws.Range("E2").Value = 30

# Size column E like the other data columns (~10.09 characters wide).
$ws.Columns.Item(5).ColumnWidth = 9.3

# The saved file leaves the cursor parked just below the new data.
$ws.Range("E3").Select() | Out-Null
